{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the manuscript revision: updates the submission date and the\n// \"tracks kept\" counts in the eps filtering-results paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Exact text replacements to perform, run by run. Using exact (non-wildcard)\n// matches keeps this robust regardless of how the paragraph's runs are split.\n// The submission date is stored as five separate runs (\"juni\" / \" \" / \"11,\" /\n// \" \" / \"2022\"); replace each run's text individually so the run structure\n// (and any per-run formatting) stays untouched, matching how the document\n// was actually revised.\nconst replacements = [\n  [\"juni\", \"January\"],\n  [\"11,\", \"20,\"],\n  [\"2022\", \"2023\"],\n  [\"eps = 50: 158 mismatches with 165 tracks kept\", \"eps = 50: 158 mismatches with 171 tracks kept\"],\n  [\"eps = 100: 103 mismatches with 147 tracks kept\", \"eps = 100: 103 mismatches with 153 tracks kept\"],\n  [\"eps = 150: 71 mismatches with 124 tracks kept\", \"eps = 150: 71 mismatches with 130 tracks kept\"],\n  [\"eps = 200: 47 mismatches with 100 tracks kept\", \"eps = 200: 47 mismatches with 106 tracks kept\"],\n  [\"eps = 250: 30 mismatches with 69 tracks kept\", \"eps = 250: 30 mismatches with 75 tracks kept\"],\n  [\"eps = 300: 9 mismatches with 49 tracks kept\", \"eps = 300: 9 mismatches with 51 tracks kept\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the manuscript revision: updates the submission date and the\n# \"tracks kept\" counts in the eps filtering-results paragraph.\n\n$d = $word.ActiveDocument\n\n# Pairs of (FindText, ReplaceWith). The submission date is stored as five\n# separate runs (\"juni\" / \" \" / \"11,\" / \" \" / \"2022\"); matching each token\n# individually (rather than the whole \"juni 11, 2022\" phrase at once) keeps\n# the replacement confined within each existing run instead of merging them\n# into one, which mirrors how the document was actually revised.\n$replacements = @(\n  @(\"juni\", \"January\"),\n  @(\"11,\", \"20,\"),\n  @(\"2022\", \"2023\"),\n  @(\"eps = 50: 158 mismatches with 165 tracks kept\", \"eps = 50: 158 mismatches with 171 tracks kept\"),\n  @(\"eps = 100: 103 mismatches with 147 tracks kept\", \"eps = 100: 103 mismatches with 153 tracks kept\"),\n  @(\"eps = 150: 71 mismatches with 124 tracks kept\", \"eps = 150: 71 mismatches with 130 tracks kept\"),\n  @(\"eps = 200: 47 mismatches with 100 tracks kept\", \"eps = 200: 47 mismatches with 106 tracks kept\"),\n  @(\"eps = 250: 30 mismatches with 69 tracks kept\", \"eps = 250: 30 mismatches with 75 tracks kept\"),\n  @(\"eps = 300: 9 mismatches with 49 tracks kept\", \"eps = 300: 9 mismatches with 51 tracks kept\")\n)\n\nforeach ($pair in $replacements) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Execute(\n    $findText,\n    $true,\n    $true,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $replaceText,\n    2\n  )\n}\n\n$d.Save()\n"}
